$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 43; existing rows 43:51 shift down to 44:52
$ws.Rows("43:43").Insert()

# Populate the newly inserted row 43 with the new record
$ws.Range("A43").Value = 5
$ws.Range("B43").Value = "Macroferia Regional de Talca"
$ws.Range("C43").Value = "Maule"
$ws.Range("D43").Value = 44505
$ws.Range("E43").Value = 7
$ws.Range("F43").Value = 100112026
$ws.Range("G43").Value = "Haba"
$ws.Range("H43").Value = "Sin especificar"
$ws.Range("I43").Value = "Primera"
$ws.Range("J43").Value = 600
$ws.Range("K43").Value = 5000
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = 5000
$ws.Range("N43").Value = "$/saco 25 kilos"
$ws.Range("O43").Value = "Región del Maule"
$ws.Range("P43").Value = 200
$ws.Range("Q43").Value = 25
$ws.Range("R43").Value = "Hortaliza"
